$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GlobalConstantIntTable")

$ws.Range("A24").Value = "TranscendGoldOne"
$ws.Range("B24").Value = 10000

$ws.Range("A25").Value = "TranscendGoldTwo"
$ws.Range("B25").Value = 20000

$ws.Range("A26").Value = "TranscendGoldThree"
$ws.Range("B26").Value = 30000

$ws.Range("A24").Select()
